$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Columns A/B use the default "General" style, so assigning a numeric-looking
    # string would otherwise get silently coerced to a number by Excel. Force
    # the cell to Text just long enough to type the value in verbatim, then put
    # the cell style back exactly the way it was (no explicit style, like the
    # original file) so styles.xml stays untouched.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 4 - Bird / Temperature / CG: Q, pval, PdeltaAIC-cov updated stats
Set-TextValue $ws.Range("A4") " 103.5"
Set-TextValue $ws.Range("B4") "0.374"
$ws.Range("F4").Value = "0.034"

# Row 6 - Bird / Temperature / CZG: updated stats
Set-TextValue $ws.Range("A6") "  81.7"
Set-TextValue $ws.Range("B6") "0.044"
$ws.Range("F6").Value = "0.396"

# Row 10 - Mammal / Temperature / CG: updated stats (now includes PdeltaAIC as a covariate)
Set-TextValue $ws.Range("A10") "   3.8"
Set-TextValue $ws.Range("B10") "0.450"
$ws.Range("F10").Value = "0.431"

# Row 13 - Mammal / Temperature / TotalCG: pval tweak
$ws.Range("F13").Value = "0.850"
